# Journal de travail update: add new log entries (rows 24-30), push the
# blank buffer rows down (31-36), and relocate the "Type / Temps total"
# summary table (Table663) + Total row from rows 30-41 to rows 41-52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New work-log entries in the Tableau2 data area (rows 24-30)
# ---------------------------------------------------------------------

$ws.Range("B24:B30").NumberFormat = "d-mmm"

$ws.Range("B24").Value = "5/17/2019"
$ws.Range("C24").Value = "Conception"
$ws.Range("D24").Value = "Correction des dernières erreurs des changements de pages."
$ws.Range("E24").Value = 1.5
$ws.Rows.Item(24).RowHeight = 30

$ws.Range("B25").Value = "5/17/2019"
$ws.Range("C25").Value = "Réalisation"
$ws.Range("D25").Value = "Création de la page d'accueil et finalisation de la page événement. Cette tâche n'aurait pas dû être effectuée avant les pages d'inscirption et de connexion. Cependant, cela était nécessaire à ma bonne compréhension de mon modèle MVC, les tâches ont donc été inversées. Le retard que je semble avoir pris est bien moins important qu'il ne paraît."
$ws.Range("E25").Value = 1.5
$ws.Rows.Item(25).RowHeight = 120

$ws.Range("B26").Value = "5/17/2019"
$ws.Range("C26").Value = "Réalisation"
$ws.Range("D26").Value = "Création de la page d'inscription."
$ws.Range("E26").Value = 0.75

$ws.Range("B27").Value = "5/21/2019"
$ws.Range("C27").Value = "Autres"
$ws.Range("D27").Value = "Rassemblement des livrables oubliés et rédaction d'un mail à destination de mes experts ainsi que de mon chef de projet pour m'excuser de mon retard dans l'envoi de mes livrables."
$ws.Range("E27").Value = 0.5
$ws.Rows.Item(27).RowHeight = 60

$ws.Range("B28").Value = "5/21/2019"
$ws.Range("C28").Value = "Réalisation"
$ws.Range("D28").Value = "Ajout sur la page d'inscription d'une fonction allant chercher la liste des écoles dans ma base de données afin de limiter le choix d'entrée de l'utilisateur. Dû à des bugs et des problèmes de fonction, la réalisation de cette petite fonctionnalité m'a pris beaucoup de temps."
$ws.Range("E28").Value = 3.25
$ws.Rows.Item(28).RowHeight = 90

$ws.Range("B29").Value = "5/21/2019"
$ws.Range("C29").Value = "Documentation"
$ws.Range("D29").Value = "Continuation de la documentation."
$ws.Range("E29").Value = 1.5

$ws.Range("B30").Value = "5/21/2019"
$ws.Range("C30").Value = "Réalisation"
$ws.Range("D30").Value = "Finalisation de la page d'inscription"
$ws.Range("E30").Value = 1.5

# ---------------------------------------------------------------------
# 2) New blank buffer rows (31-36), matching the two that used to sit
#    at 24-25 right under the data.
# ---------------------------------------------------------------------

$ws.Range("B31:B34").NumberFormat = "d-mmm"

# ---------------------------------------------------------------------
# 3) Move the "Type / Temps total" summary table + Total row down from
#    rows 30-41 to rows 41-52.
# ---------------------------------------------------------------------

$ws.Range("C30:D41").Clear()

$ws.Range("C41").Value = "Type"
$ws.Range("D41").Value = "Temps total"

$ws.Range("C42").Value = "Analyse"
$ws.Range("D42").Formula = "=SUMIF(Tableau2[Type d'activité],C42,Tableau2[Temps nécessaire])"

$ws.Range("C43").Value = "Recherches"
$ws.Range("D43").Formula = "=SUMIF(Tableau2[Type d'activité],C43,Tableau2[Temps nécessaire])"

$ws.Range("C44").Value = "Conception"
$ws.Range("D44").Formula = "=SUMIF(Tableau2[Type d'activité],C44,Tableau2[Temps nécessaire])"

$ws.Range("C45").Value = "Documentation"
$ws.Range("D45").Formula = "=SUMIF(Tableau2[Type d'activité],C45,Tableau2[Temps nécessaire])"

$ws.Range("C46").Value = "Réalisation"
$ws.Range("D46").Formula = "=SUMIF(Tableau2[Type d'activité],C46,Tableau2[Temps nécessaire])"

$ws.Range("C47").Value = "Autres"
$ws.Range("D47").Formula = "=SUMIF(Tableau2[Type d'activité],C47,Tableau2[Temps nécessaire])"

$ws.Range("C48").Value = "Tests"
$ws.Range("D48").Formula = "=SUMIF(Tableau2[Type d'activité],C48,Tableau2[Temps nécessaire])"

$ws.Range("C49").Value = "Planification"
$ws.Range("D49").Formula = "=SUMIF(Tableau2[Type d'activité],C49,Tableau2[Temps nécessaire])"

$ws.Range("C50").Value = "Absences"
$ws.Range("D50").Formula = "=SUMIF(Tableau2[Type d'activité],C50,Tableau2[Temps nécessaire])"

$ws.Range("C52").Value = "Total"
$ws.Range("C52").Font.Bold = $true
$ws.Range("D52").Formula = "=SUM(Table663[Temps total])"
$ws.Range("D52").Font.Bold = $true

# ---------------------------------------------------------------------
# 4) Resize the two tables (ListObjects) to their new extents.
# ---------------------------------------------------------------------

$tblLog = $ws.ListObjects.Item("Tableau2")
$tblLog.Resize($ws.Range("B3:E36"))

$tblSummary = $ws.ListObjects.Item("Table663")
$tblSummary.Resize($ws.Range("C41:D50"))

# ---------------------------------------------------------------------
# 5) Data validation on the "Type d'activité" column now spans the
#    bigger data range and points at the relocated summary list.
# ---------------------------------------------------------------------

$ws.Range("C4:C25").Validation.Delete()
$ws.Range("C4:C36").Validation.Add(3, 1, 1, "=`$C`$42:`$C`$50")

# ---------------------------------------------------------------------
# 6) Selection / scroll position, matching where the author was working.
# ---------------------------------------------------------------------

$ws.Rows.Item(31).Select()
